$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''69.164.25'
$ws.Range('E2').Value = '  +0.96%  '

# Row 3
$ws.Range('D3').Value = '''3.771.98'
$ws.Range('E3').Value = '  -1.07%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = '''631.20'
$ws.Range('E5').Value = '  +2.91%  '

# Row 6
$ws.Range('D6').Value = '''166.83'
$ws.Range('E6').Value = '  +2.08%  '

# Row 7
$ws.Range('D7').Value = '''3.768.73'
$ws.Range('E7').Value = '  -1.15%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('E9').Value = '  +0.75%  '

# Row 10
$ws.Range('E10').Value = '  -0.85%  '

# Row 11
$ws.Range('E11').Value = '  +2.17%  '

# Row 12
$ws.Range('D12').Value = '''6.75'
$ws.Range('E12').Value = '  -0.86%  '

# Row 13
$ws.Range('E13').Value = '  -3.33%  '

# Row 14
$ws.Range('D14').Value = '''35.07'
$ws.Range('E14').Value = '  +0.04%  '

# Row 15
$ws.Range('D15').Value = '''4.409.70'
$ws.Range('E15').Value = '  -0.84%  '

# Row 16
$ws.Range('D16').Value = '''3.782.31'
$ws.Range('E16').Value = '  -0.19%  '

# Row 17
$ws.Range('D17').Value = '''69.193.96'
$ws.Range('E17').Value = '  +1.07%  '

# Row 18
$ws.Range('D18').Value = '''17.60'
$ws.Range('E18').Value = '  -2.59%  '

# Row 19
$ws.Range('E19').Value = '  +0.26%  '

# Row 20
$ws.Range('D20').Value = '''7.01'
$ws.Range('E20').Value = '  -0.70%  '

# Row 21
$ws.Range('D21').Value = '''462.83'
$ws.Range('E21').Value = '  +0.01%  '

# Row 22
$ws.Range('E22').Value = '  -0.93%  '

# Row 23
$ws.Range('D23').Value = '''0.705'
$ws.Range('E23').Value = '  +1.14%  '

# Row 24
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '''0.0000145'
$ws.Range('E24').Value = '  -1.24%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''82.51'
$ws.Range('E25').Value = '  -1.00%  '

# Row 26
$ws.Range('D26').Value = '''12.08'
$ws.Range('E26').Value = '  +1.05%  '

# Row 27
$ws.Range('D27').Value = '''2.14'
$ws.Range('E27').Value = '  +2.08%  '

# Row 28
$ws.Range('E28').Value = '  +1.15%  '

# Row 29
$ws.Range('E29').Value = '  -0.11%  '

# Row 30
$ws.Range('D30').Value = '''3.923.23'

# Row 31
$ws.Range('D31').Value = '''2.31'
$ws.Range('E31').Value = '  +5.28%  '

# Row 32
$ws.Range('E32').Value = '  +2.79%  '

# Row 33
$ws.Range('E33').Value = '  -1.71%  '

# Row 34
$ws.Range('D34').Value = '''0.178'
$ws.Range('E34').Value = '  +22.91%  '

# Row 35
$ws.Range('E35').Value = '  -1.37%  '

# Row 36
$ws.Range('E36').Value = '  -0.03%  '

# Row 37
$ws.Range('E37').Value = '  -0.80%  '

# Row 38
$ws.Range('E38').Value = '  -1.06%  '

# Row 39
$ws.Range('E39').Value = '  +0.80%  '

# Row 40
$ws.Range('E40').Value = '  +6.95%  '

# Row 41
$ws.Range('E41').Value = '  -1.25%  '

# Row 42
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.960'
$ws.Range('E42').Value = '  -2.23%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.20%  '

# Row 45
$ws.Range('D45').Value = '''158.15'
$ws.Range('E45').Value = '  +3.23%  '

# Row 46
$ws.Range('E46').Value = '  +7.35%  '

# Row 47
$ws.Range('E47').Value = '  +2.80%  '

# Row 48
$ws.Range('D48').Value = '''43.42'
$ws.Range('E48').Value = '  +1.23%  '

# Row 49
$ws.Range('D49').Value = '''47.12'
$ws.Range('E49').Value = '  +1.10%  '

# Row 50
$ws.Range('E50').Value = '  +0.22%  '

# Row 51
$ws.Range('D51').Value = '''8.37'
$ws.Range('E51').Value = '  +0.23%  '
